$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cd14"
$ws.Cells.Item(2, 3).Value = "Itgb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 54.82987066666667
$ws.Cells.Item(2, 8).Value = 164.489612
$ws.Cells.Item(2, 9).Value = 0.9762630652055824
$ws.Cells.Item(2, 10).Value = 0.9762630652055824
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 112.513392
$ws.Cells.Item(2, 14).Value = 337.540176
$ws.Cells.Item(2, 15).Value = 0.3275312977368564
$ws.Cells.Item(2, 16).Value = 0.3275312977368564
$ws.Cells.Item(2, 17).Value = 6169.094731627968
$ws.Cells.Item(2, 18).Value = 55521.85258465172
$ws.Cells.Item(2, 19).Value = 0.3197567086793456
$ws.Cells.Item(2, 20).Value = 0.3197567086793456

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cd14"
$ws.Cells.Item(3, 3).Value = "Itgb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 54.82987066666667
$ws.Cells.Item(3, 8).Value = 164.489612
$ws.Cells.Item(3, 9).Value = 0.9762630652055824
$ws.Cells.Item(3, 10).Value = 0.9762630652055824
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 106.314466
$ws.Cells.Item(3, 14).Value = 318.943398
$ws.Cells.Item(3, 15).Value = 0.3094859589441663
$ws.Cells.Item(3, 16).Value = 0.3094859589441664
$ws.Cells.Item(3, 17).Value = 5829.208420775731
$ws.Cells.Item(3, 18).Value = 52462.87578698158
$ws.Cells.Item(3, 19).Value = 0.3021397109169208
$ws.Cells.Item(3, 20).Value = 0.3021397109169209

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Cd14"
$ws.Cells.Item(4, 3).Value = "Itgb1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 54.82987066666667
$ws.Cells.Item(4, 8).Value = 164.489612
$ws.Cells.Item(4, 9).Value = 0.9762630652055824
$ws.Cells.Item(4, 10).Value = 0.9762630652055824
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 124.6916553333333
$ws.Cells.Item(4, 14).Value = 374.074966
$ws.Cells.Item(4, 15).Value = 0.3629827433189773
$ws.Cells.Item(4, 16).Value = 0.3629827433189773
$ws.Cells.Item(4, 17).Value = 6836.827335139244
$ws.Cells.Item(4, 18).Value = 61531.44601625321
$ws.Cells.Item(4, 19).Value = 0.3543666456093159
$ws.Cells.Item(4, 20).Value = 0.354366645609316

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cd14"
$ws.Cells.Item(5, 3).Value = "Itgb1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.9891043333333333
$ws.Cells.Item(5, 8).Value = 2.967313
$ws.Cells.Item(5, 9).Value = 0.01761131325912771
$ws.Cells.Item(5, 10).Value = 0.01761131325912771
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 112.513392
$ws.Cells.Item(5, 14).Value = 337.540176
$ws.Cells.Item(5, 15).Value = 0.3275312977368564
$ws.Cells.Item(5, 16).Value = 0.3275312977368564
$ws.Cells.Item(5, 17).Value = 111.287483585232
$ws.Cells.Item(5, 18).Value = 1001.587352267088
$ws.Cells.Item(5, 19).Value = 0.005768256286612402
$ws.Cells.Item(5, 20).Value = 0.005768256286612402

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Cd14"
$ws.Cells.Item(6, 3).Value = "Itgb1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.9891043333333333
$ws.Cells.Item(6, 8).Value = 2.967313
$ws.Cells.Item(6, 9).Value = 0.01761131325912771
$ws.Cells.Item(6, 10).Value = 0.01761131325912771
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 106.314466
$ws.Cells.Item(6, 14).Value = 318.943398
$ws.Cells.Item(6, 15).Value = 0.3094859589441663
$ws.Cells.Item(6, 16).Value = 0.3094859589441664
$ws.Cells.Item(6, 17).Value = 105.1560990166193
$ws.Cells.Item(6, 18).Value = 946.4048911495739
$ws.Cells.Item(6, 19).Value = 0.005450454172267249
$ws.Cells.Item(6, 20).Value = 0.005450454172267249

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Cd14"
$ws.Cells.Item(7, 3).Value = "Itgb1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.9891043333333333
$ws.Cells.Item(7, 8).Value = 2.967313
$ws.Cells.Item(7, 9).Value = 0.01761131325912771
$ws.Cells.Item(7, 10).Value = 0.01761131325912771
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 124.6916553333333
$ws.Cells.Item(7, 14).Value = 374.074966
$ws.Cells.Item(7, 15).Value = 0.3629827433189773
$ws.Cells.Item(7, 16).Value = 0.3629827433189773
$ws.Cells.Item(7, 17).Value = 123.3330566207064
$ws.Cells.Item(7, 18).Value = 1109.997509586358
$ws.Cells.Item(7, 19).Value = 0.006392602800248053
$ws.Cells.Item(7, 20).Value = 0.006392602800248054

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Cd14"
$ws.Cells.Item(8, 3).Value = "Itgb1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.3440333333333334
$ws.Cells.Item(8, 8).Value = 1.0321
$ws.Cells.Item(8, 9).Value = 0.00612562153528991
$ws.Cells.Item(8, 10).Value = 0.006125621535289909
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 112.513392
$ws.Cells.Item(8, 14).Value = 337.540176
$ws.Cells.Item(8, 15).Value = 0.3275312977368564
$ws.Cells.Item(8, 16).Value = 0.3275312977368564
$ws.Cells.Item(8, 17).Value = 38.7083572944
$ws.Cells.Item(8, 18).Value = 348.3752156496
$ws.Cells.Item(8, 19).Value = 0.002006332770898339
$ws.Cells.Item(8, 20).Value = 0.002006332770898338

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Cd14"
$ws.Cells.Item(9, 3).Value = "Itgb1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.3440333333333334
$ws.Cells.Item(9, 8).Value = 1.0321
$ws.Cells.Item(9, 9).Value = 0.00612562153528991
$ws.Cells.Item(9, 10).Value = 0.006125621535289909
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 106.314466
$ws.Cells.Item(9, 14).Value = 318.943398
$ws.Cells.Item(9, 15).Value = 0.3094859589441663
$ws.Cells.Item(9, 16).Value = 0.3094859589441664
$ws.Cells.Item(9, 17).Value = 36.57572011953334
$ws.Cells.Item(9, 18).Value = 329.1814810758
$ws.Cells.Item(9, 19).Value = 0.001895793854978234
$ws.Cells.Item(9, 20).Value = 0.001895793854978234

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Cd14"
$ws.Cells.Item(10, 3).Value = "Itgb1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.3440333333333334
$ws.Cells.Item(10, 8).Value = 1.0321
$ws.Cells.Item(10, 9).Value = 0.00612562153528991
$ws.Cells.Item(10, 10).Value = 0.006125621535289909
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 124.6916553333333
$ws.Cells.Item(10, 14).Value = 374.074966
$ws.Cells.Item(10, 15).Value = 0.3629827433189773
$ws.Cells.Item(10, 16).Value = 0.3629827433189773
$ws.Cells.Item(10, 17).Value = 42.89808582317779
$ws.Cells.Item(10, 18).Value = 386.0827724086
$ws.Cells.Item(10, 19).Value = 0.002223494909413337
$ws.Cells.Item(10, 20).Value = 0.002223494909413337
